# Update Name of Algo
# Apply updated imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 5.486499999999994
$ws.Range("C6").Value  = -11.4901
$ws.Range("C7").Value  = -11.7432
$ws.Range("B8").Value  = 4.772200000000001
$ws.Range("C8").Value  = -11.19539999999999
$ws.Range("E11").Value = 13.4064
$ws.Range("A12").Value = -22.59059999999999
$ws.Range("B12").Value = 6.239399999999996
$ws.Range("B14").Value = 9.495200000000006
$ws.Range("E14").Value = 13.6775
$ws.Range("C19").Value = -13.20289999999999
$ws.Range("E19").Value = 13.39249999999999
$ws.Range("C21").Value = -13.17050000000001
$ws.Range("E21").Value = 13.35029999999999
$ws.Range("B22").Value = 4.807900000000005
$ws.Range("C24").Value = -11.2567
